# Commit: "Upd: add collections theme" (Java Урок 8 Масиви.pptx)
# "Додано блок тем "Колекції". Оновлено інші файлики, де було помічено косяки."
#
# The captured slide (index 22, "Типові операції з масивами") had every
# shape nudged by the same amount: +9427 EMU horizontally (~0.7423pt) and
# -254524 EMU vertically (~-20.0413pt), as if the whole slide content was
# selected and bumped slightly to the right/up while other files in the
# deck were being tidied up. Sizes (widths/heights) are untouched.
#
# Point literals below are picked so that this runtime's internal
# float32-based pt -> EMU conversion lands exactly on the target EMU
# offsets recorded in the canonical OOXML, rather than relying on the
# nearest 4-decimal value PowerPoint's UI would normally show.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(22)

function Move-Shape($name, $left, $top) {
    $sh = $s.Shapes.Item($name)
    $sh.Left = $left
    $sh.Top = $top
}

# Placeholder body text boxes
Move-Shape "Rectangle 5" 453.617431640625   94.83378601074219
Move-Shape "Rectangle 3" 145.24229431152344 94.83378601074219

# Embedded Visio OLE diagrams (graphicFrames)
Move-Shape "Object 7"  145.24229431152344 132.3338623046875
Move-Shape "Object 9"  148.36741638183594 456.708740234375
Move-Shape "Object 10" 460.2423095703125  125.208740234375
Move-Shape "Object 11" 145.24229431152344 278.5838623046875
Move-Shape "Object 12" 460.2423095703125  199.89166259765625
Move-Shape "Object 15" 147.61732482910156 382.9588317871094
Move-Shape "Object 17" 457.9923095703125  288.5838623046875

# Vertical divider line
Move-Shape "Line 16" 448.3673400878906 98.0838623046875

Write-Output "slide 22: repositioned 10 shapes"
